$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Егор"
$ws.Range("C2").Value = "Губин"
$ws.Range("D2").Value = "Вячеславович"
$ws.Range("E2").Value = "fadey"
$ws.Range("F2").Value = 123
$ws.Range("G2").Value = "komrad.gubi2017@yandex.ru"

$ws.Range("B3").Value = "Илья"
$ws.Range("C3").Value = "Сибелев"
$ws.Range("D3").Value = "Владимирович"
$ws.Range("E3").Value = "ghost"
$ws.Range("F3").Value = 123
$ws.Range("G3").Value = "sibelev@yandex.ru"

$ws.Range("B4").Value = "Кира"
$ws.Range("C4").Value = "Амеличева"
$ws.Range("D4").Value = "Александровна"
$ws.Range("E4").Value = "kaa"
$ws.Range("F4").Value = "'111"
$ws.Range("G4").Value = "kaa@mail.ru"

$ws.Rows("5:7").Delete()
